$d = $word.ActiveDocument

$replacements = @(
    @("354÷7=", "360÷3="),
    @("416÷6=", "868÷5="),
    @("887÷2=", "234÷9="),
    @("444÷3=", "315÷2="),
    @("898÷2=", "334÷4="),
    @("861÷4=", "559÷3="),
    @("539÷7=", "225÷2="),
    @("994÷5=", "981÷4="),
    @("901÷8=", "746÷5="),
    @("760÷9=", "600÷6="),
    @("359÷8=", "484÷3="),
    @("461÷7=", "762÷3="),
    @("192÷8=", "996÷3="),
    @("541÷3=", "450÷3="),
    @("966÷5=", "894÷9="),
    @("620÷6=", "779÷7="),
    @("330÷7=", "405÷7="),
    @("513÷7=", "515÷7="),
    @("237÷2=", "100÷7="),
    @("223÷5=", "638÷2="),
    @("740÷3=", "242÷7="),
    @("401÷3=", "926÷6="),
    @("118÷4=", "821÷4="),
    @("943÷2=", "273÷2="),
    @("696÷9=", "314÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
